$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, avoiding Excel's automatic
# number/date reinterpretation of numeric-looking strings (e.g. "1.003").
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "22.368.34"
Set-TextValue "E2" "  -0.34%  "

# Row 3
Set-TextValue "D3" "1.566.43"
Set-TextValue "E3" "  -0.48%  "

# Row 4
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.28%  "

# Row 5
Set-TextValue "D5" "1.003"
Set-TextValue "E5" "  +0.25%  "

# Row 6
Set-TextValue "D6" "290.83"
Set-TextValue "E6" "  +0.29%  "

# Row 7
Set-TextValue "D7" "0.3773"
Set-TextValue "E7" "  +2.50%  "

# Row 8
Set-TextValue "D8" "49.22"
Set-TextValue "E8" "  -0.09%  "

# Row 9
Set-TextValue "D9" "0.3402"

# Row 10
Set-TextValue "D10" "0.07589"
Set-TextValue "E10" "  -0.97%  "

# Row 11
Set-TextValue "D11" "1.138"
Set-TextValue "E11" "  -3.10%  "

# Row 12
Set-TextValue "D12" "1.003"
Set-TextValue "E12" "  +0.26%  "

# Row 13
Set-TextValue "D13" "21.03"
Set-TextValue "E13" "  -1.43%  "

# Row 14
Set-TextValue "D14" "5.981"
Set-TextValue "E14" "  -1.54%  "

# Row 15
Set-TextValue "D15" "6.906"
Set-TextValue "E15" "  -0.39%  "

# Row 16
Set-TextValue "D16" "1.567.04"
Set-TextValue "E16" "  -0.46%  "

# Row 17
Set-TextValue "D17" "0.00001133"
Set-TextValue "E17" "  -0.11%  "

# Row 18
Set-TextValue "D18" "89.88"
Set-TextValue "E18" "  -0.48%  "

# Row 19
Set-TextValue "D19" "0.06742"
Set-TextValue "E19" "  +0.08%  "

# Row 20
Set-TextValue "E20" "  +0.28%  "

# Row 21
Set-TextValue "E21" "  +0.24%  "

# Row 22
Set-TextValue "D22" "6.202"
Set-TextValue "E22" "  -1.24%  "

# Row 23
Set-TextValue "D23" "11.94"
Set-TextValue "E23" "  -0.91%  "

# Row 24
Set-TextValue "D24" "22.358.26"
Set-TextValue "E24" "  -0.40%  "

# Row 25
Set-TextValue "D25" "2.399"
Set-TextValue "E25" "  +1.07%  "

# Row 26
Set-TextValue "D26" "2.691"
Set-TextValue "E26" "  -6.72%  "

# Row 27
Set-TextValue "E27" "  +0.10%  "

# Row 28
Set-TextValue "D28" "147.57"
Set-TextValue "E28" "  +0.31%  "

# Row 29
Set-TextValue "D29" "5.022"
Set-TextValue "E29" "  +0.65%  "

# Row 30
Set-TextValue "D30" "125.90"
Set-TextValue "E30" "  -0.06%  "

# Row 31
Set-TextValue "D31" "1.739.81"
Set-TextValue "E31" "  -0.48%  "

# Row 32
Set-TextValue "D32" "2.013"
Set-TextValue "E32" "  -0.28%  "

# Row 33
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.9905"
Set-TextValue "E33" "  -3.66%  "

# Row 34
Set-TextValue "B34" "Filecoin"
Set-TextValue "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D34" "6.047"
Set-TextValue "E34" "  -3.32%  "

# Row 35
Set-TextValue "D35" "10.08"
Set-TextValue "E35" "  -0.45%  "

# Row 36
Set-TextValue "D36" "1.436"
Set-TextValue "E36" "  +10.37%  "

# Row 37
Set-TextValue "E37" "  -0.11%  "

# Row 38
Set-TextValue "D38" "0.02503"
Set-TextValue "E38" "  -1.65%  "

# Row 39
Set-TextValue "E39" "  -1.46%  "

# Row 40
Set-TextValue "D40" "0.06461"
Set-TextValue "E40" "  -0.45%  "

# Row 41
Set-TextValue "D41" "5.400"
Set-TextValue "E41" "  -2.73%  "

# Row 42
Set-TextValue "D42" "0.6314"
Set-TextValue "E42" "  -0.98%  "

# Row 43
Set-TextValue "D43" "11.30"
Set-TextValue "E43" "  -3.78%  "

# Row 44
Set-TextValue "E44" "  +0.21%  "

# Row 45
Set-TextValue "D45" "13.99"
Set-TextValue "E45" "  -2.40%  "

# Row 46
Set-TextValue "D46" "3.805"
Set-TextValue "E46" "  +1.16%  "

# Row 47
Set-TextValue "D47" "0.5932"
Set-TextValue "E47" "  -1.17%  "

# Row 48
Set-TextValue "D48" "2.079"
Set-TextValue "E48" "  -1.68%  "

# Row 49
Set-TextValue "D49" "1.255"
Set-TextValue "E49" "  -0.35%  "

# Row 50
Set-TextValue "D50" "124.60"
Set-TextValue "E50" "  -0.47%  "

# Row 51
Set-TextValue "D51" "0.07320"
Set-TextValue "E51" "  +0.26%  "
